# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates the DAMSLTag (column I) and DialogAct (column J) values for the rows
# whose dialog-act annotations changed after the transcript clean up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 16; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 23; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 26; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 30; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 33; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 55; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 57; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 73; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 80; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 88; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 94; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 109; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 110; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 128; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 129; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 132; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 135; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 137; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 143; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 147; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 157; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 159; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 163; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 167; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 168; I = 'qy'; J = 'Yes-No-Question' },
    @{ Row = 174; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 176; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 179; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 194; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 200; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 201; I = 'qy'; J = 'Yes-No-Question' },
    @{ Row = 209; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 211; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 227; I = '%'; J = 'Uninterpretable' },
    @{ Row = 231; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 232; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 248; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 249; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 257; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 280; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 283; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 305; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 310; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 312; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 313; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 315; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 319; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 333; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 334; I = '%'; J = 'Uninterpretable' },
    @{ Row = 336; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 348; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 351; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 357; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 365; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 377; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 387; I = 'ba'; J = 'Appreciation' },
    @{ Row = 392; I = 'ba'; J = 'Appreciation' },
    @{ Row = 404; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 410; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 414; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 423; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 432; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 440; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 445; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 467; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 468; I = 'sd'; J = 'Statement-non-opinion' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
